# Lift state updating : progress, ButtonCallbackData to enum
# Appends 4 new lift rows (IDs 24-27) to the "lifts" worksheet, following
# the same NONE/S/O/N pattern used by the previous rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lifts")

$startRow = 27
$startId = 24
$rowsToAdd = 4

for ($i = 0; $i -lt $rowsToAdd; $i++) {
    $row = $startRow + $i
    $id = $startId + $i

    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 2).Value = "NONE"
    $ws.Cells.Item($row, 3).Value = "S"
    $ws.Cells.Item($row, 4).Value = "O"
    $ws.Cells.Item($row, 5).Value = "N"
}
